$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert the Q:V probability columns for data rows 2-41 back to their
# previous values (this undoes the earlier "fix" that shifted them).
$newQ = 0.98750000000000004
$newR = 0.25750000000000001
$newS = 0.017499999999999998
$newT = 0.017499999999999998
$newU = 0.017499999999999998
$newV = 0.48749999999999999

for ($row = 2; $row -le 41; $row++) {
    $ws.Range("Q$row").Value = $newQ
    $ws.Range("R$row").Value = $newR
    $ws.Range("S$row").Value = $newS
    $ws.Range("T$row").Value = $newT
    $ws.Range("U$row").Value = $newU
    $ws.Range("V$row").Value = $newV
}

# Update the active selection / view on Sheet1 to match the reverted commit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N8").Select()
